# CRM_1 sheet1: refresh the sample login data.
# Old:  A1 = jesmi   B1 = admin
# New:  A1 = jesmi1  B1 = admin  C1 = CHINNU
# Shared-string insertion order below matches the target workbook's table
# (admin, CHINNU, jesmi1) so indices line up with the expected diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "admin"
$ws.Range("C1").Value = "CHINNU"
$ws.Range("A1").Value = "jesmi1"
